$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "58.719.92"
$ws.Range("E2").Value = "  +2.05%  "

# Row 3
$ws.Range("D3").Value = "3.087.61"
$ws.Range("E3").Value = "  +0.16%  "

# Row 4
$ws.Range("E4").Value = "  -0.01%  "

# Row 5
$ws.Range("D5").Value = "'523.04"
$ws.Range("E5").Value = "  +1.25%  "

# Row 6
$ws.Range("D6").Value = "'143.67"
$ws.Range("E6").Value = "  +0.44%  "

# Row 7
$ws.Range("E7").Value = "  +0.00%  "

# Row 8
$ws.Range("D8").Value = "'0.439"
$ws.Range("E8").Value = "  +0.92%  "

# Row 9
$ws.Range("D9").Value = "'7.35"
$ws.Range("E9").Value = "  +0.90%  "

# Row 10
$ws.Range("E10").Value = "  +0.64%  "

# Row 11
$ws.Range("E11").Value = "  +2.87%  "

# Row 12
$ws.Range("D12").Value = "3.618.63"
$ws.Range("E12").Value = "  +0.10%  "

# Row 13
$ws.Range("E13").Value = "  +1.15%  "

# Row 14
$ws.Range("D14").Value = "'26.80"
$ws.Range("E14").Value = "  +3.92%  "

# Row 15
$ws.Range("D15").Value = "'0.0000167"
$ws.Range("E15").Value = "  +0.92%  "

# Row 16
$ws.Range("D16").Value = "58.732.35"
$ws.Range("E16").Value = "  +1.83%  "

# Row 17
$ws.Range("D17").Value = "3.086.84"
$ws.Range("E17").Value = "  -0.11%  "

# Row 18
$ws.Range("E18").Value = "  +0.55%  "

# Row 19
$ws.Range("E19").Value = "  -1.29%  "

# Row 20
$ws.Range("D20").Value = "'8.14"
$ws.Range("E20").Value = "  -0.50%  "

# Row 21
$ws.Range("D21").Value = "'342.31"
$ws.Range("E21").Value = "  +1.32%  "

# Row 22
$ws.Range("E22").Value = "  -0.04%  "

# Row 23
$ws.Range("D23").Value = "'0.506"
$ws.Range("E23").Value = "  +0.67%  "

# Row 24
$ws.Range("D24").Value = "'65.76"
$ws.Range("E24").Value = "  +0.31%  "

# Row 25
$ws.Range("E25").Value = "  -0.40%  "

# Row 26
$ws.Range("D26").Value = "'1.00"
$ws.Range("E26").Value = "  -0.22%  "

# Row 27
$ws.Range("E27").Value = "  -0.45%  "

# Row 28
$ws.Range("E28").Value = "  +2.44%  "

# Row 29
$ws.Range("E29").Value = "  +1.85%  "

# Row 30
$ws.Range("D30").Value = "'1.85"
$ws.Range("E30").Value = "  +1.83%  "

# Row 31
$ws.Range("D31").Value = "'21.04"
$ws.Range("E31").Value = "  +0.90%  "

# Row 32
$ws.Range("E32").Value = "  +2.71%  "

# Row 33
$ws.Range("D33").Value = "'154.30"
$ws.Range("E33").Value = "  +0.18%  "

# Row 34
$ws.Range("D34").Value = "'4.61"
$ws.Range("E34").Value = "  +2.02%  "

# Row 35
$ws.Range("D35").Value = "'6.13"
$ws.Range("E35").Value = "  +3.89%  "

# Row 36
$ws.Range("D36").Value = "'26.91"
$ws.Range("E36").Value = "  -3.43%  "

# Row 37
$ws.Range("D37").Value = "'1.30"
$ws.Range("E37").Value = "  +5.20%  "

# Row 38
$ws.Range("D38").Value = "'0.0682"
$ws.Range("E38").Value = "  -1.01%  "

# Row 39
$ws.Range("D39").Value = "3.130.09"
$ws.Range("E39").Value = "  +0.11%  "

# Row 40
$ws.Range("E40").Value = "  +1.14%  "

# Row 41
$ws.Range("D41").Value = "'36.80"
$ws.Range("E41").Value = "  -0.16%  "

# Row 42
$ws.Range("B42").Value = "Stacks"
$ws.Range("C42").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D42").Value = "'1.49"
$ws.Range("E42").Value = "  +7.35%  "

# Row 43
$ws.Range("B43").Value = "FirstDigitalUSD"
$ws.Range("C43").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D43").Value = "'1.00"
$ws.Range("E43").Value = "  -0.03%  "

# Row 44
$ws.Range("B44").Value = "Mantle"
$ws.Range("C44").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D44").Value = "'0.665"
$ws.Range("E44").Value = "  -0.81%  "

# Row 45
$ws.Range("D45").Value = "2.272.32"
$ws.Range("E45").Value = "  -0.66%  "

# Row 46
$ws.Range("D46").Value = "'0.0256"
$ws.Range("E46").Value = "  +1.61%  "

# Row 47
$ws.Range("D47").Value = "'20.77"
$ws.Range("E47").Value = "  +2.13%  "

# Row 48
$ws.Range("D48").Value = "'0.964"
$ws.Range("E48").Value = "  +1.80%  "

# Row 49
$ws.Range("D49").Value = "'6.02"
$ws.Range("E49").Value = "  +2.16%  "

# Row 50
$ws.Range("D50").Value = "'0.749"
$ws.Range("E50").Value = "  +8.06%  "

# Row 51
$ws.Range("D51").Value = "'265.50"
$ws.Range("E51").Value = "  +11.46%  "
